# Update "Glosario violencia de género" workbook:
# The description for "Violencia sexual" (row 6, column B) had a manual
# line break removed (it now reads as a single paragraph), which also
# shrinks the wrapped row height from 3 lines to 2 lines.
# We also move the current selection to B7 to match where the author
# left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("Violencia sexual"): remove the embedded line break between
# "pleno" and "de la sexualidad de las mujeres" so the text reads as one
# continuous sentence.
$ws.Range("B6").Value = "La violencia sexual se manifiesta de diversas formas: física, simbólica, explícita o implícitamente, y corresponde a cualquier práctica que atente contra el desarrollo pleno de la sexualidad de las mujeres. Como expresión del continuo de violencia, está presente desde la infancia y se manifiesta en diversos espacios: familia, instituciones educativas, religiosas, laborales, espacios públicos, centros de salud, etc."

# The row no longer needs to be as tall now that the manual break is gone.
$ws.Rows(6).RowHeight = 47.25

# Move the selection to B7, as left by the author after editing B6, and
# scroll the view down so row 4 is the first fully visible row.
$ws.Range("B7").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
